# Update scripts with new TPM values.
# The underlying NATMI ligand/receptor metrics (ligand & receptor average/total
# expression, derived specificities, and edge weights/specificities) were
# recomputed against the updated TPM data. This patches the already-computed
# values in the LR-pairs worksheet to reflect the refreshed numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.120168333333333
$ws.Range("H2").Value = 3.360505
$ws.Range("I2").Value = 0.001768092629909379
$ws.Range("J2").Value = 0.001768092629909379
$ws.Range("M2").Value = 91.51130433333333
$ws.Range("N2").Value = 274.533913
$ws.Range("O2").Value = 0.9685519820468944
$ws.Range("P2").Value = 0.9685519820468945
$ws.Range("Q2").Value = 102.5080652562294
$ws.Range("R2").Value = 922.5725873060649
$ws.Range("S2").Value = 0.001712489621141235
$ws.Range("T2").Value = 0.001712489621141235
$ws.Range("G3").Value = 1.120168333333333
$ws.Range("H3").Value = 3.360505
$ws.Range("I3").Value = 0.001768092629909379
$ws.Range("J3").Value = 0.001768092629909379
$ws.Range("O3").Value = 0.001425786415744213
$ws.Range("P3").Value = 0.001425786415744214
$ws.Range("Q3").Value = 0.15090011652
$ws.Range("R3").Value = 1.35810104868
$ws.Range("S3").Value = 0.000002520922453502253
$ws.Range("T3").Value = 0.000002520922453502253
$ws.Range("G4").Value = 1.120168333333333
$ws.Range("H4").Value = 3.360505
$ws.Range("I4").Value = 0.001768092629909379
$ws.Range("J4").Value = 0.001768092629909379
$ws.Range("M4").Value = 2.836578333333333
$ws.Range("N4").Value = 8.509734999999999
$ws.Range("O4").Value = 0.03002223153736139
$ws.Range("P4").Value = 0.03002223153736139
$ws.Range("Q4").Value = 3.177445224019444
$ws.Range("R4").Value = 28.597007016175
$ws.Range("S4").Value = 0.0000530820863146416
$ws.Range("T4").Value = 0.0000530820863146416
$ws.Range("I5").Value = 0.9534130698726969
$ws.Range("J5").Value = 0.9534130698726969
$ws.Range("M5").Value = 91.51130433333333
$ws.Range("N5").Value = 274.533913
$ws.Range("O5").Value = 0.9685519820468944
$ws.Range("P5").Value = 0.9685519820468945
$ws.Range("Q5").Value = 55275.68382413404
$ws.Range("R5").Value = 497481.1544172064
$ws.Range("S5").Value = 0.9234301185346149
$ws.Range("T5").Value = 0.923430118534615
$ws.Range("I6").Value = 0.9534130698726969
$ws.Range("J6").Value = 0.9534130698726969
$ws.Range("O6").Value = 0.001425786415744213
$ws.Range("P6").Value = 0.001425786415744214
$ws.Range("S6").Value = 0.00135936340361748
$ws.Range("T6").Value = 0.00135936340361748
$ws.Range("I7").Value = 0.9534130698726969
$ws.Range("J7").Value = 0.9534130698726969
$ws.Range("M7").Value = 2.836578333333333
$ws.Range("N7").Value = 8.509734999999999
$ws.Range("O7").Value = 0.03002223153736139
$ws.Range("P7").Value = 0.03002223153736139
$ws.Range("Q7").Value = 1713.381841052064
$ws.Range("S7").Value = 0.02862358793446463
$ws.Range("T7").Value = 0.02862358793446463
$ws.Range("H8").Value = 85.18441
$ws.Range("I8").Value = 0.04481883749739363
$ws.Range("J8").Value = 0.04481883749739363
$ws.Range("M8").Value = 91.51130433333333
$ws.Range("N8").Value = 274.533913
$ws.Range("O8").Value = 0.9685519820468944
$ws.Range("P8").Value = 0.9685519820468945
$ws.Range("Q8").Value = 2598.445489321814
$ws.Range("R8").Value = 23386.00940389633
$ws.Range("S8").Value = 0.04340937389113827
$ws.Range("T8").Value = 0.04340937389113828
$ws.Range("H9").Value = 85.18441
$ws.Range("I9").Value = 0.04481883749739363
$ws.Range("J9").Value = 0.04481883749739363
$ws.Range("O9").Value = 0.001425786415744213
$ws.Range("P9").Value = 0.001425786415744214
$ws.Range("Q9").Value = 3.82512074664
$ws.Range("R9").Value = 34.42608671976
$ws.Range("S9").Value = 0.00006390208967323121
$ws.Range("T9").Value = 0.00006390208967323121
$ws.Range("H10").Value = 85.18441
$ws.Range("I10").Value = 0.04481883749739363
$ws.Range("J10").Value = 0.04481883749739363
$ws.Range("M10").Value = 2.836578333333333
$ws.Range("N10").Value = 8.509734999999999
$ws.Range("O10").Value = 0.03002223153736139
$ws.Range("P10").Value = 0.03002223153736139
$ws.Range("S10").Value = 0.001345561516582126
$ws.Range("T10").Value = 0.001345561516582126
